$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '65.796.91'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.52%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.677.05'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.66%  '

$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '600.61'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.79%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '156.62'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.29%  '

$ws.Range("E7").Value = '  +0.04%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.587'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.66%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.123'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.17%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '5.90'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.08%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.397'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -3.92%  '

$ws.Range("E12").Value = '  +0.14%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '29.45'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.84%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000201'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +5.15%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.162.46'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.75%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.634.85'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.62%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.694.90'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.30%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.59'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.54%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.82'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.29%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.61'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.89%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '350.95'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.89%  '

$ws.Range("E22").Value = '  -0.10%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '70.33'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.17%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.0000110'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +5.81%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.82'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.85%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.62'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -6.19%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.170'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.32%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.62'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.55%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.16'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.93%  '

$ws.Range("E30").Value = '  +0.03%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.17'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.00%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '531.31'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -4.20%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.76'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -4.50%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.56'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.66%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.41'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -5.06%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.426'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.91%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '20.47'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.01%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '160.73'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.81%  '

$ws.Range("B39").Value = 'FirstDigitalUSD'
$ws.Range("C39").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.00'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.01%  '

$ws.Range("B40").Value = 'Stacks'
$ws.Range("C40").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.96'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.84%  '

$ws.Range("E41").Value = '  +0.04%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '42.24'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.87%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '165.95'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.32%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.09'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.14%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0621'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.51%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '23.07'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.03%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.24'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.88%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0263'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.07%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.651'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.94%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '20.36'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.41%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0984'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.49%  '
